$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 47. This shifts the existing rows 47..79
# down to 48..80, preserving all of their original data untouched.
$ws.Rows(47).Insert()

# Populate the newly inserted row 47 with the new weekly record.
$ws.Cells.Item(47, 1).Value  = 10
$ws.Cells.Item(47, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(47, 3).Value  = "La Araucanía"
$ws.Cells.Item(47, 4).Value  = 45001
$ws.Cells.Item(47, 5).Value  = 9
$ws.Cells.Item(47, 6).Value  = "Fruta"
$ws.Cells.Item(47, 7).Value  = 100107
$ws.Cells.Item(47, 8).Value  = "Otros"
$ws.Cells.Item(47, 9).Value  = 100107011
$ws.Cells.Item(47, 10).Value = "Tuna"
$ws.Cells.Item(47, 11).Value = "Sin especificar"
$ws.Cells.Item(47, 12).Value = "Primera"
$ws.Cells.Item(47, 13).Value = 300
$ws.Cells.Item(47, 14).Value = 16000
$ws.Cells.Item(47, 15).Value = 18000
$ws.Cells.Item(47, 16).Value = 17000
$ws.Cells.Item(47, 17).Value = "$/caja 16 kilos"
$ws.Cells.Item(47, 18).Value = "Provincia de Los Andes"
$ws.Cells.Item(47, 19).Value = 1062
$ws.Cells.Item(47, 20).Value = 16
